$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 1933.3334
$ws.Range("I69").Value = 900
$ws.Range("J69").Value = 4000
$ws.Range("K69").Value = 2700
$ws.Range("L69").Value = 12000
$ws.Range("M69").Value = -1826
$ws.Range("N69").Value = -13748

$ws.Range("H72").Value = 1933.3334
$ws.Range("I72").Value = 900
$ws.Range("J72").Value = 4000
$ws.Range("K72").Value = 8100
$ws.Range("L72").Value = 36000
$ws.Range("M72").Value = -3732
$ws.Range("N72").Value = -44736

$ws.Range("H96").Value = 381.7143
$ws.Range("I96").Value = 273.83334
$ws.Range("J96").Value = 1029
$ws.Range("K96").Value = 821.5000200000001
$ws.Range("L96").Value = 3087
$ws.Range("M96").Value = 551.4999799999999
$ws.Range("N96").Value = -5833

$ws.Range("H111").Value = 2051.5
$ws.Range("I111").Value = 2051.5
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 6154.5
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -3087.5

$ws.Range("H112").Value = 2119.8125
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2119.8125
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 6359.4375
$ws.Range("N112").Value = -8575.4375

$ws.Range("H132").Value = 6672939.5
$ws.Range("I132").Value = 7755839.5
$ws.Range("J132").Value = 20840.715
$ws.Range("K132").Value = 23267518.5
$ws.Range("L132").Value = 62522.145
$ws.Range("M132").Value = -23264988.5
$ws.Range("N132").Value = -67582.145

$ws.Range("H137").Value = 1488.5483
$ws.Range("I137").Value = 1382.5883
$ws.Range("J137").Value = 1617.2142
$ws.Range("K137").Value = 4147.7649
$ws.Range("L137").Value = 4851.642599999999
$ws.Range("M137").Value = -1597.7649
$ws.Range("N137").Value = -9951.642599999999

$ws.Range("H138").Value = 3019.4546
$ws.Range("I138").Value = 3109.8333
$ws.Range("J138").Value = 2999.3704
$ws.Range("K138").Value = 9329.499899999999
$ws.Range("L138").Value = 8998.111199999999
$ws.Range("M138").Value = -4189.499899999999
$ws.Range("N138").Value = -19278.1112

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14654.797
$ws.Range("I32").Value = 13458.366
$ws.Range("J32").Value = 15470.546
$ws.Range("K32").Value = 13458.366
$ws.Range("L32").Value = 15470.546
$ws.Range("M32").Value = -13171.366

$ws.Range("H132").Value = 2351.2327
$ws.Range("I132").Value = 1988.0303
$ws.Range("J132").Value = 3549.8
$ws.Range("K132").Value = 5964.090899999999
$ws.Range("L132").Value = 10649.4
$ws.Range("M132").Value = -3434.090899999999
$ws.Range("N132").Value = -15709.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0

$ws.Range("H80").Value = 537.94446
$ws.Range("I80").Value = 193
$ws.Range("J80").Value = 813.9
$ws.Range("K80").Value = 193
$ws.Range("L80").Value = 813.9
$ws.Range("M80").Value = 805

$ws.Range("H83").Value = 537.94446
$ws.Range("I83").Value = 193
$ws.Range("J83").Value = 813.9
$ws.Range("K83").Value = 965
$ws.Range("L83").Value = 4069.5
$ws.Range("M83").Value = 4027

$ws.Range("H94").Value = 12501072
$ws.Range("I94").Value = 14706891
$ws.Range("J94").Value = 1433.3334
$ws.Range("K94").Value = 14706891
$ws.Range("L94").Value = 1433.3334
$ws.Range("M94").Value = -14706440

$ws.Range("H105").Value = 77685230
$ws.Range("I105").Value = 91809460
$ws.Range("J105").Value = 1995
$ws.Range("K105").Value = 91809460
$ws.Range("L105").Value = 1995
$ws.Range("M105").Value = -91807713
$ws.Range("N105").Value = -5489

$ws.Range("H134").Value = 5091.2173
$ws.Range("I134").Value = 777.1818
$ws.Range("J134").Value = 100000
$ws.Range("K134").Value = 2331.5454
$ws.Range("L134").Value = 300000
$ws.Range("M134").Value = 203.4546
$ws.Range("N134").Value = -305070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 258.9091
$ws.Range("I7").Value = 218.5
$ws.Range("J7").Value = 366.66666
$ws.Range("K7").Value = 218.5
$ws.Range("L7").Value = 366.66666
$ws.Range("M7").Value = -105.5
$ws.Range("N7").Value = -592.66666

$ws.Range("H62").Value = 22224478
$ws.Range("I62").Value = 2580
$ws.Range("J62").Value = 50001850
$ws.Range("K62").Value = 2580
$ws.Range("L62").Value = 50001850
$ws.Range("M62").Value = -1956
$ws.Range("N62").Value = -50003098

$ws.Range("H65").Value = 22224478
$ws.Range("I65").Value = 2580
$ws.Range("J65").Value = 50001850
$ws.Range("K65").Value = 12900
$ws.Range("L65").Value = 250009250
$ws.Range("M65").Value = -9780
$ws.Range("N65").Value = -250015490

$ws.Range("H86").Value = 2801334.2
$ws.Range("I86").Value = 4779272
$ws.Range("J86").Value = 32221.2
$ws.Range("K86").Value = 4779272
$ws.Range("L86").Value = 32221.2
$ws.Range("M86").Value = -4778149
$ws.Range("N86").Value = -34467.2

$ws.Range("H89").Value = 2801334.2
$ws.Range("I89").Value = 4779272
$ws.Range("J89").Value = 32221.2
$ws.Range("K89").Value = 23896360
$ws.Range("L89").Value = 161106
$ws.Range("M89").Value = -23890744
$ws.Range("N89").Value = -172338

$ws.Range("H105").Value = 368.7353
$ws.Range("I105").Value = 319.96295
$ws.Range("J105").Value = 556.8570999999999
$ws.Range("K105").Value = 319.96295
$ws.Range("L105").Value = 556.8570999999999
$ws.Range("M105").Value = 1427.03705

$ws.Range("H107").Value = 1541
$ws.Range("I107").Value = 867.5
$ws.Range("J107").Value = 1990
$ws.Range("K107").Value = 867.5
$ws.Range("L107").Value = 1990
$ws.Range("M107").Value = 1052.5

$ws.Range("H132").Value = 1994.439
$ws.Range("I132").Value = 1652.909
$ws.Range("J132").Value = 3403.25
$ws.Range("K132").Value = 4958.727000000001
$ws.Range("L132").Value = 10209.75
$ws.Range("M132").Value = -2428.727000000001
$ws.Range("N132").Value = -15269.75

$ws.Range("H134").Value = 8773165
$ws.Range("I134").Value = 1141.4878
$ws.Range("J134").Value = 31251474
$ws.Range("K134").Value = 3424.463400000001
$ws.Range("L134").Value = 93754422
$ws.Range("M134").Value = -889.4634000000005

$ws.Range("H135").Value = 38276.668
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 38276.668
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 38276.668
$ws.Range("N135").Value = -48416.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 261.27777
$ws.Range("I14").Value = 261.27777
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 783.83331
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -610.83331

$ws.Range("H97").Value = 837.3333
$ws.Range("I97").Value = 723
$ws.Range("J97").Value = 894.5
$ws.Range("K97").Value = 2169
$ws.Range("L97").Value = 2683.5
$ws.Range("M97").Value = -1673
$ws.Range("N97").Value = -3675.5

$ws.Range("H113").Value = 689.1
$ws.Range("I113").Value = 593.5333000000001
$ws.Range("J113").Value = 746.4400000000001
$ws.Range("K113").Value = 1780.5999
$ws.Range("L113").Value = 2239.32
$ws.Range("M113").Value = 389.4000999999998
$ws.Range("N113").Value = -6579.32

$ws.Range("H122").Value = 899.0192
$ws.Range("I122").Value = 665.4286
$ws.Range("J122").Value = 985.0789
$ws.Range("K122").Value = 5988.8574
$ws.Range("L122").Value = 8865.7101
$ws.Range("M122").Value = -3538.8574
$ws.Range("N122").Value = -13765.7101

$ws.Range("H131").Value = 23846248
$ws.Range("I131").Value = 250000370
$ws.Range("J131").Value = 40551.973
$ws.Range("K131").Value = 750001110
$ws.Range("L131").Value = 121655.919
$ws.Range("M131").Value = -749996070
$ws.Range("N131").Value = -131735.919

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2373.9736
$ws.Range("I102").Value = 1583.826
$ws.Range("J102").Value = 3585.5334
$ws.Range("K102").Value = 1583.826
$ws.Range("L102").Value = 3585.5334
$ws.Range("M102").Value = 38.17399999999998
$ws.Range("N102").Value = -6829.5334

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0

$ws.Range("H132").Value = 11007.267
$ws.Range("I132").Value = 15442.556
$ws.Range("J132").Value = 4354.3335
$ws.Range("K132").Value = 46327.66800000001
$ws.Range("L132").Value = 13063.0005
$ws.Range("M132").Value = -43797.66800000001
$ws.Range("N132").Value = -18123.0005

$ws.Range("H138").Value = 38520
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 38520
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 38520
$ws.Range("N138").Value = -48800

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1700.7778
$ws.Range("I22").Value = 2950.5
$ws.Range("J22").Value = 1343.7142
$ws.Range("K22").Value = 2950.5
$ws.Range("L22").Value = 1343.7142
$ws.Range("M22").Value = -2655.5
$ws.Range("N22").Value = -1933.7142

$ws.Range("H27").Value = 1700.7778
$ws.Range("I27").Value = 2950.5
$ws.Range("J27").Value = 1343.7142
$ws.Range("K27").Value = 2950.5
$ws.Range("L27").Value = 1343.7142
$ws.Range("M27").Value = -2843.5
$ws.Range("N27").Value = -1557.7142

$ws.Range("H46").Value = 7999.2856
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 7999.2856
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 7999.2856
$ws.Range("N46").Value = -8375.285599999999

$ws.Range("H74").Value = 29217
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 29217
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 29217
$ws.Range("N74").Value = -31213

$ws.Range("H77").Value = 29217
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 29217
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 87651
$ws.Range("N77").Value = -97635

$ws.Range("H122").Value = 141666670
$ws.Range("I122").Value = 141666670
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 425000010
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -424997560

$ws.Range("H132").Value = 2734.1333
$ws.Range("I132").Value = 2431.4
$ws.Range("J132").Value = 3339.6
$ws.Range("K132").Value = 7294.200000000001
$ws.Range("L132").Value = 10018.8
$ws.Range("M132").Value = -4764.200000000001
$ws.Range("N132").Value = -15078.8

$ws.Range("H136").Value = 1870.8518
$ws.Range("I136").Value = 1503.7059
$ws.Range("J136").Value = 2495
$ws.Range("K136").Value = 4511.1177
$ws.Range("L136").Value = 7485
$ws.Range("M136").Value = -1961.1177
$ws.Range("N136").Value = -12585

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 352.66666
$ws.Range("I100").Value = 359.25
$ws.Range("J100").Value = 300
$ws.Range("K100").Value = 718.5
$ws.Range("L100").Value = 600
$ws.Range("M100").Value = -177.5
$ws.Range("N100").Value = -1682

$ws.Range("H126").Value = 92593790
$ws.Range("I126").Value = 111112350
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 333337050
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -333334580

$ws.Range("H136").Value = 1025.9667
$ws.Range("I136").Value = 741
$ws.Range("J136").Value = 1518.1818
$ws.Range("K136").Value = 2223
$ws.Range("L136").Value = 4554.5454
$ws.Range("M136").Value = 327
$ws.Range("N136").Value = -9654.545399999999

# Clear cells that should be removed
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N111").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N61").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
